$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Yearly")

# Row 64 - British Airways: A 967 -> 968, D 137 -> 138
$ws.Range("A64").Value = 968
$ws.Range("D64").Value = 138

# Row 65 - Delta Air Lines: A 978 -> 979
$ws.Range("A65").Value = 979

# Row 66 - EasyJet: A 980 -> 981
$ws.Range("A66").Value = 981

# Row 67 - Jet2: A 1007 -> 1009
$ws.Range("A67").Value = 1009

# Row 68 - KLM: A 1010 -> 1012
$ws.Range("A68").Value = 1012

# Row 69 - Lufthansa: A 1022 -> 1024
$ws.Range("A69").Value = 1024

# Row 70 - Ryanair: A 1048 -> 1050
$ws.Range("A70").Value = 1050

# Row 71 - United Airlines: A 1071 -> 1074
$ws.Range("A71").Value = 1074
